$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural: insert two new columns (A becomes new, and a new column
# appears between the old "gmail" and "pass" columns) ---
$ws.Columns("A").Insert()
$ws.Columns("C").Insert()

# --- Row 1 headers: new "yahoo" and "name" columns ---
$ws.Range("A1").Value = "yahoo"
$ws.Range("C1").Value = "name"
$ws.Range("I1").Value = "access token"
$ws.Range("J1").Value = "access token secret"

# --- Row 2: fill "name" column with same value as username column ---
$ws.Range("C2").Value = "goddard0001"

# --- New row 3: second yahoo/gmail-style bot account ---
$ws.Range("C3").Value = "doodlebob0042"
$ws.Range("D3").Value = "dgeg45?12"

# Copy the date-format (style 1) from E2 down into E3, but leave it blank
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)

$ws.Range("G3").Value = "wRVSTKYDJWCHFmagYg8qhZQwF"
$ws.Range("I3").Value = "890052028675784000aHiBYsLVMFbyWxzvU1oEndDY5Pjx641"
$ws.Range("J3").Value = "ZSWgsMuQRFY1zzPjB8HNAWyWFzrBqD69gN2E687NJsP90"

# Copy the "value" style (style 2: Arial 8 FF333333) onto the plain new cells
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("J3").PasteSpecial(-4122)

# H3 gets the same base style plus word-wrap + vertical-centering
$ws.Range("H3").Value = "zGSRwPz9w3SpPW1XWoyurJMxU1Rqk9vg1vD0dOVw64N8u278Ss"
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").WrapText = $true
$ws.Range("H3").VerticalAlignment = -4108
$ws.Rows(3).RowHeight = 20.4

# --- New row 4: a third account name, reusing row-3 password ---
$ws.Range("A4").Value = "doodlebob0043"
$ws.Range("D4").Value = "dgeg45?12"

# --- New row 6: another bot entry with its own styled cells ---
$ws.Range("C6").Value = "cougarGuy01"
$ws.Range("D6").Value = "meow01"
$ws.Range("G6").Value = "MEpbDNplsX724pFmaKtTjlMFW"
$ws.Range("H6").Value = "BEIyPrCYzWe8X65fArKE65cXAkLtyUd8RAYlLllHrw4ubiND9T"
$ws.Range("I6").Value = "932834480406122497-kEdX6h2yZ7ocWhEgVznG22My7qbeLx5"
$ws.Range("J6").Value = "0hjm4BHFfCiHgNwNVccNw66puL0ZN22hwGinrtdjrW2d6"
$ws.Range("K6").Value = "application name = cougarGuy01!"
$ws.Range("L6").Value = "url = https://na.com"

$ws.Range("G2").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("J6").PasteSpecial(-4122)

# B6 (email) gets the Arial style too, but in the lighter gray color
$ws.Range("B6").Value = "bavalley1@cougars.ccis.edu"
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Font.Color = 6710886

# Column C ("name") needs an explicit width (copied visually from column B)
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Tidy up: clear clipboard marching ants, select the last-edited cell
$excel.CutCopyMode = $false
$ws.Range("J6").Select()
